# Atualização automática de JULIO_DE_CASTILHOS.xlsx
#
# The commit removes the "Desarquivamentos Pendentes" sheet entirely and
# renames two other sheets to their upper-case equivalents:
#   "Paineis DARQ"              -> "PAINEIS DARQ"
#   "Recolhimento x Eliminacao" -> "RECOLHIMENTO X ELIMINAÇÃO"

$wb = $excel.ActiveWorkbook

# 1) Delete the sheet that is no longer part of the workbook.
$wb.Worksheets("Desarquivamentos Pendentes").Delete()

# 2) Rename the remaining sheets (content/format of these sheets is
#    otherwise untouched).
$wb.Worksheets("Paineis DARQ").Name = "PAINEIS DARQ"
$wb.Worksheets("Recolhimento x Eliminacao").Name = "RECOLHIMENTO X ELIMINAÇÃO"

# 3) Keep the originally-active sheet selected/active (the first tab),
#    since deleting the last sheet in the workbook otherwise shifts the
#    active tab to the new last sheet.
$wb.Worksheets("PAINEIS DARQ").Activate()
